$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 85.833336
$ws.Range("I8").Value = 85.833336
$ws.Range("K8").Value = 257.500008
$ws.Range("M8").Value = -118.500008
$ws.Range("H99").Value = 578.1818
$ws.Range("I99").Value = 311.42856
$ws.Range("J99").Value = 1045
$ws.Range("K99").Value = 934.28568
$ws.Range("L99").Value = 3135
$ws.Range("M99").Value = 563.71432
$ws.Range("N99").Value = -6131
$ws.Range("H129").Value = 1387.25
$ws.Range("I129").Value = 293.33334
$ws.Range("J129").Value = 1543.5238
$ws.Range("K129").Value = 880.0000200000001
$ws.Range("L129").Value = 4630.5714
$ws.Range("M129").Value = 4119.99998
$ws.Range("N129").Value = -14630.5714
$ws.Range("H135").Value = 1279.2709
$ws.Range("I135").Value = 312.46875
$ws.Range("J135").Value = 3212.875
$ws.Range("K135").Value = 2812.21875
$ws.Range("L135").Value = 28915.875
$ws.Range("M135").Value = -277.21875
$ws.Range("N135").Value = -33985.875
$ws.Range("H138").Value = 3987.2983
$ws.Range("I138").Value = 6609.591
$ws.Range("J138").Value = 2339
$ws.Range("K138").Value = 19828.773
$ws.Range("L138").Value = 7017
$ws.Range("M138").Value = -14688.773
$ws.Range("N138").Value = -17297

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 882
$ws.Range("I45").Value = 882
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 882
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -505
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 2829.1304
$ws.Range("I61").Value = 1940.7142
$ws.Range("J61").Value = 4211.1113
$ws.Range("K61").Value = 1940.7142
$ws.Range("L61").Value = 4211.1113
$ws.Range("M61").Value = -1728.7142
$ws.Range("N61").Value = -4635.1113
$ws.Range("H74").Value = 4654
$ws.Range("I74").Value = 877
$ws.Range("J74").Value = 10834.546
$ws.Range("K74").Value = 877
$ws.Range("L74").Value = 10834.546
$ws.Range("M74").Value = -3
$ws.Range("N74").Value = -12582.546
$ws.Range("H77").Value = 4654
$ws.Range("I77").Value = 877
$ws.Range("J77").Value = 10834.546
$ws.Range("K77").Value = 4385
$ws.Range("L77").Value = 54172.73
$ws.Range("M77").Value = -17
$ws.Range("N77").Value = -62908.73
$ws.Range("H132").Value = 1906.303
$ws.Range("I132").Value = 1882
$ws.Range("J132").Value = 1971.1111
$ws.Range("K132").Value = 5646
$ws.Range("L132").Value = 5913.3333
$ws.Range("M132").Value = -3116
$ws.Range("N132").Value = -10973.3333
$ws.Range("H136").Value = 2829.1304
$ws.Range("I136").Value = 1940.7142
$ws.Range("J136").Value = 4211.1113
$ws.Range("K136").Value = 5822.142599999999
$ws.Range("L136").Value = 12633.3339
$ws.Range("M136").Value = -3272.142599999999
$ws.Range("N136").Value = -17733.3339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6255
$ws.Range("I105").Value = 6255
$ws.Range("K105").Value = 6255
$ws.Range("M105").Value = -4508
$ws.Range("H134").Value = 92244.73
$ws.Range("I134").Value = 101359.2
$ws.Range("K134").Value = 304077.6
$ws.Range("M134").Value = -301542.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3700
$ws.Range("I16").Value = 3980
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 3980
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -3693
$ws.Range("N16").Value = -3574
$ws.Range("H58").Value = 2286.9583
$ws.Range("I58").Value = 1617.6364
$ws.Range("K58").Value = 1617.6364
$ws.Range("M58").Value = -1414.6364
$ws.Range("H105").Value = 854.5454999999999
$ws.Range("I105").Value = 811.6667
$ws.Range("K105").Value = 811.6667
$ws.Range("M105").Value = 935.3333
$ws.Range("H113").Value = 3700
$ws.Range("I113").Value = 3980
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 3980
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -1810
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 4515.125
$ws.Range("I132").Value = 4046.4
$ws.Range("K132").Value = 12139.2
$ws.Range("M132").Value = -9609.200000000001
$ws.Range("H134").Value = 4186.8
$ws.Range("I134").Value = 5380.2
$ws.Range("K134").Value = 16140.6
$ws.Range("M134").Value = -13605.6
$ws.Range("H136").Value = 2286.9583
$ws.Range("I136").Value = 1617.6364
$ws.Range("K136").Value = 4852.9092
$ws.Range("M136").Value = -2302.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3525.6
$ws.Range("J94").Value = 5026.6665
$ws.Range("L94").Value = 15079.9995
$ws.Range("N94").Value = -16431.9995
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H117").Value = 700
$ws.Range("I117").Value = 700
$ws.Range("K117").Value = 2100
$ws.Range("M117").Value = 1342
$ws.Range("H133").Value = 8333.333000000001
$ws.Range("I133").Value = 5000
$ws.Range("J133").Value = 10000
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 30000
$ws.Range("M133").Value = -9940
$ws.Range("N133").Value = -40120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 368.33334
$ws.Range("I13").Value = 302.5
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 302.5
$ws.Range("L13").Value = 500
$ws.Range("M13").Value = -163.5
$ws.Range("N13").Value = -778
$ws.Range("H113").Value = 1411.0555
$ws.Range("I113").Value = 904.3
$ws.Range("J113").Value = 2044.5
$ws.Range("K113").Value = 904.3
$ws.Range("L113").Value = 2044.5
$ws.Range("M113").Value = 1265.7
$ws.Range("N113").Value = -6384.5
$ws.Range("H132").Value = 2429.25
$ws.Range("I132").Value = 1607.3334
$ws.Range("J132").Value = 2922.4
$ws.Range("K132").Value = 4822.0002
$ws.Range("L132").Value = 8767.200000000001
$ws.Range("M132").Value = -2292.0002
$ws.Range("N132").Value = -13827.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1867.2084
$ws.Range("I16").Value = 2064.1904
$ws.Range("J16").Value = 488.33334
$ws.Range("K16").Value = 2064.1904
$ws.Range("L16").Value = 488.33334
$ws.Range("M16").Value = -1894.1904
$ws.Range("N16").Value = -828.33334
$ws.Range("H40").Value = 1947.4
$ws.Range("I40").Value = 1947.4
$ws.Range("K40").Value = 1947.4
$ws.Range("M40").Value = -1811.4
$ws.Range("H68").Value = 1698.9656
$ws.Range("I68").Value = 1679.5238
$ws.Range("J68").Value = 1750
$ws.Range("K68").Value = 1679.5238
$ws.Range("L68").Value = 1750
$ws.Range("M68").Value = -930.5237999999999
$ws.Range("N68").Value = -3248
$ws.Range("H70").Value = 33000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 33000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 33000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -33540
$ws.Range("H71").Value = 1698.9656
$ws.Range("I71").Value = 1679.5238
$ws.Range("J71").Value = 1750
$ws.Range("K71").Value = 8397.618999999999
$ws.Range("L71").Value = 8750
$ws.Range("M71").Value = -4653.618999999999
$ws.Range("N71").Value = -16238
$ws.Range("H73").Value = 33000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 33000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 33000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -34872
$ws.Range("H119").Value = 21204
$ws.Range("J119").Value = 21204
$ws.Range("L119").Value = 21204
$ws.Range("N119").Value = -30880
$ws.Range("H136").Value = 1082.8536
$ws.Range("I136").Value = 955.62067
$ws.Range("K136").Value = 2866.86201
$ws.Range("M136").Value = -316.8620099999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 16267375
$ws.Range("J2").Value = 18356500
$ws.Range("L2").Value = 18356500
$ws.Range("N2").Value = -18356724
$ws.Range("H12").Value = 35219
$ws.Range("I12").Value = 5150
$ws.Range("J12").Value = 50253.5
$ws.Range("K12").Value = 5150
$ws.Range("L12").Value = 50253.5
$ws.Range("M12").Value = -5008
$ws.Range("N12").Value = -50537.5
$ws.Range("H132").Value = 2380.25
$ws.Range("I132").Value = 2435.0688
$ws.Range("J132").Value = 2274.2666
$ws.Range("K132").Value = 7305.2064
$ws.Range("L132").Value = 6822.7998
$ws.Range("M132").Value = -4775.2064
$ws.Range("N132").Value = -11882.7998
$ws.Range("H136").Value = 6950.1763
$ws.Range("I136").Value = 7613.3667
$ws.Range("K136").Value = 22840.1001
$ws.Range("M136").Value = -20290.1001
